$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the Price column (D) are plain text in the source data (e.g.
# '42.338.22' or '20.71'), even though some look numeric. A leading apostrophe
# forces Excel to store the assigned value as literal text (matching the
# original inline-string cells) instead of silently re-interpreting it as a
# number (which would corrupt values like multi-dot prices or drop trailing
# zeros, e.g. '20.70' -> 20.7).

$ws.Range("D2").Value = '''42.373.25'
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").Value = '''2.252.67'
$ws.Range("E3").Value = '  -1.15%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '''247.75'
$ws.Range("E5").Value = '  -1.35%  '

$ws.Range("E6").Value = '  -3.51%  '

$ws.Range("D7").Value = '''74.21'
$ws.Range("E7").Value = '  -1.32%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -4.13%  '

$ws.Range("D10").Value = '''42.02'
$ws.Range("E10").Value = '  +6.40%  '

$ws.Range("D11").Value = '''0.0945'
$ws.Range("E11").Value = '  -3.76%  '

$ws.Range("D12").Value = '''7.16'
$ws.Range("E12").Value = '  -3.89%  '

$ws.Range("E13").Value = '  -3.50%  '

$ws.Range("D14").Value = '''2.586.32'
$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("E15").Value = '  -3.97%  '

$ws.Range("D16").Value = '''0.856'
$ws.Range("E16").Value = '  -1.31%  '

$ws.Range("D17").Value = '''2.245.53'
$ws.Range("E17").Value = '  -1.49%  '

$ws.Range("D18").Value = '''42.177.01'
$ws.Range("E18").Value = '  -1.14%  '

$ws.Range("E19").Value = '  -2.04%  '

$ws.Range("D20").Value = '''6.12'
$ws.Range("E20").Value = '  -1.59%  '

$ws.Range("E21").Value = '  -0.84%  '

$ws.Range("D22").Value = '''2.27'
$ws.Range("E22").Value = '  +5.08%  '

$ws.Range("D23").Value = '''230.36'
$ws.Range("E23").Value = '  -2.62%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '''8.19'
$ws.Range("E25").Value = '  +27.35%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''11.19'
$ws.Range("E26").Value = '  -0.98%  '

$ws.Range("D27").Value = '''3.58'
$ws.Range("E27").Value = '  -7.38%  '

$ws.Range("D28").Value = '''2.31'
$ws.Range("E28").Value = '  -3.45%  '

$ws.Range("D29").Value = '''2.17'
$ws.Range("E29").Value = '  +2.45%  '

$ws.Range("D30").Value = '''169.25'
$ws.Range("E30").Value = '  +1.03%  '

$ws.Range("D31").Value = '''20.70'
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("E32").Value = '  -5.75%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = '''0.119'
$ws.Range("E33").Value = '  -5.64%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '''30.68'
$ws.Range("E34").Value = '  -2.82%  '

$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("D36").Value = '''4.55'
$ws.Range("E36").Value = '  +0.57%  '

$ws.Range("E37").Value = '  +3.25%  '

$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("D39").Value = '''13.48'
$ws.Range("E39").Value = '  -1.03%  '

$ws.Range("D40").Value = '''2.19'
$ws.Range("E40").Value = '  -4.70%  '

$ws.Range("D41").Value = '''5.80'
$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("E42").Value = '  -2.70%  '

$ws.Range("D43").Value = '''61.42'
$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("D44").Value = '''107.62'
$ws.Range("E44").Value = '  +2.08%  '

$ws.Range("D45").Value = '''8.65'
$ws.Range("E45").Value = '  -3.61%  '

$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("D47").Value = '''0.998'
$ws.Range("E47").Value = '  -0.22%  '

$ws.Range("E48").Value = '  -4.22%  '

$ws.Range("E49").Value = '  -0.54%  '

$ws.Range("D50").Value = '''2.29'
$ws.Range("E50").Value = '  +0.75%  '

$ws.Range("D51").Value = '''4.12'
$ws.Range("E51").Value = '  -2.77%  '
